# Update the "scraped_at" timestamps (column K) on the "snapshot" sheet.
# These are re-scrape timestamps that changed between the previous commit
# and this one; every other cell/value in the workbook is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$updates = @(
    @{Row=2;  Value="2025-11-23T03:01:50.494957+00:00"},
    @{Row=3;  Value="2025-11-23T03:01:50.494995+00:00"},
    @{Row=4;  Value="2025-11-23T03:01:53.135597+00:00"},
    @{Row=5;  Value="2025-11-23T03:01:53.135629+00:00"},
    @{Row=6;  Value="2025-11-23T03:01:55.857638+00:00"},
    @{Row=7;  Value="2025-11-23T03:01:58.572287+00:00"},
    @{Row=8;  Value="2025-11-23T03:02:00.951950+00:00"},
    @{Row=9;  Value="2025-11-23T03:02:00.951980+00:00"},
    @{Row=10; Value="2025-11-23T03:02:00.951999+00:00"},
    @{Row=11; Value="2025-11-23T03:02:03.585113+00:00"},
    @{Row=12; Value="2025-11-23T03:02:05.879462+00:00"},
    @{Row=13; Value="2025-11-23T03:02:08.627346+00:00"},
    @{Row=14; Value="2025-11-23T03:02:10.915961+00:00"},
    @{Row=15; Value="2025-11-23T03:02:13.741583+00:00"},
    @{Row=16; Value="2025-11-23T03:02:18.786635+00:00"},
    @{Row=17; Value="2025-11-23T03:02:18.786666+00:00"},
    @{Row=18; Value="2025-11-23T03:02:21.489255+00:00"},
    @{Row=19; Value="2025-11-23T03:02:21.489289+00:00"},
    @{Row=20; Value="2025-11-23T03:02:21.489308+00:00"},
    @{Row=21; Value="2025-11-23T03:02:23.825051+00:00"},
    @{Row=22; Value="2025-11-23T03:02:23.825082+00:00"},
    @{Row=23; Value="2025-11-23T03:02:26.584481+00:00"},
    @{Row=24; Value="2025-11-23T03:02:26.584513+00:00"},
    @{Row=25; Value="2025-11-23T03:02:26.584532+00:00"},
    @{Row=26; Value="2025-11-23T03:02:26.584552+00:00"},
    @{Row=27; Value="2025-11-23T03:02:29.362414+00:00"},
    @{Row=28; Value="2025-11-23T03:02:29.362445+00:00"},
    @{Row=29; Value="2025-11-23T03:02:32.003050+00:00"},
    @{Row=30; Value="2025-11-23T03:02:32.003080+00:00"},
    @{Row=31; Value="2025-11-23T03:02:32.003099+00:00"},
    @{Row=32; Value="2025-11-23T03:02:32.003121+00:00"},
    @{Row=33; Value="2025-11-23T03:02:34.766140+00:00"},
    @{Row=34; Value="2025-11-23T03:02:34.766172+00:00"},
    @{Row=35; Value="2025-11-23T03:02:40.337064+00:00"},
    @{Row=36; Value="2025-11-23T03:02:40.337093+00:00"},
    @{Row=37; Value="2025-11-23T03:02:42.545420+00:00"},
    @{Row=38; Value="2025-11-23T03:02:42.545450+00:00"}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 11).Value = $u.Value
}
